# daily auto push: 2026-01-15 06:50 UTC
#
# A new data row for 2026/01/15 (Thursday, hour 13, ranking 28) was
# inserted at sheet row 629, pushing the existing rows 629-670
# (2026/12/29 .. 2027/01/05) down to rows 630-671. The sheet's
# dimension grows from A1:D670 to A1:D671.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (old rows 629..670) down by one row so a new
# row can be inserted at 629 without overwriting anything.
$ws.Rows("629:629").Insert()

# Force columns A/B on the new row to be read back as text (matching
# the rest of the "日付"/"曜日" columns) instead of Excel's automatic
# date/locale conversion, then write the values.
$ws.Range("A629:B629").NumberFormat = "@"
$ws.Range("A629").Value = "2026/01/15"
$ws.Range("B629").Value = "木"
$ws.Range("C629").Value = 13
$ws.Range("D629").Value = 28

# Drop the temporary text formatting again so the new row ends up
# style-less, just like every other data row in the sheet.
$ws.Range("A629:D629").ClearFormats()
